$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Respostas Mais Acertadas" answer list: drop "b-Um ser humano" and
# append a new "a-Um animal" entry after "b-Brincando na praia" (the other
# labels shift up by one row to fill the gap). ---
$ws.Range("A10").Value = "b-Carros"
$ws.Range("A11").Value = "b-Brincando na praia"
$ws.Range("A12").Value = "a-Um animal"

# --- Updated totals / counts ---
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 8

$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 2

$ws.Range("B16").Value = 5

$ws.Range("B21").Value = 4

$ws.Range("B27").Value = 4
